# Apply the "Updated symbol list on Tue Jan 31 09:59:56 UTC 2023 with GitHub Actions" edit.
# Refreshes the Coin / Link / Price / Volume(1h) columns on Sheet1 with the new scrape values.
# D/E columns hold numeric- and percent-looking text, so a leading apostrophe is used to
# force Excel to keep them as literal text (matching the original inlineStr formatting)
# instead of silently converting them into numbers and losing the exact printed precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'310.70"
$ws.Range("E2").Value = "'0.55%"

# Row 3
$ws.Range("D3").Value = "'37.15"
$ws.Range("E3").Value = "'-1.96%"

# Row 4
$ws.Range("D4").Value = "'5.128"
$ws.Range("E4").Value = "'1.01%"

# Row 5
$ws.Range("D5").Value = "'0.07773"
$ws.Range("E5").Value = "'-1.50%"

# Row 6
$ws.Range("D6").Value = "'4.375"
$ws.Range("E6").Value = "'-1.14%"

# Row 7
$ws.Range("D7").Value = "'8.223"
$ws.Range("E7").Value = "'-0.56%"

# Row 8
$ws.Range("E8").Value = "'-7.25%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9224"
$ws.Range("E9").Value = "'-1.02%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1220"
$ws.Range("E10").Value = "'-4.88%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1913"
$ws.Range("E11").Value = "'1.58%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09244"
$ws.Range("E12").Value = "'6.00%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03430"
$ws.Range("E13").Value = "'-0.79%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09676"
$ws.Range("E14").Value = "'0.22%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001370"
$ws.Range("E15").Value = "'-2.30%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006018"
$ws.Range("E16").Value = "'-3.29%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.560"
$ws.Range("E17").Value = "'-0.93%"

# Row 18
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.937"
$ws.Range("E18").Value = "'-5.70%"

# Row 19
$ws.Range("D19").Value = "'0.3405"
$ws.Range("E19").Value = "'-0.87%"

# Row 20
$ws.Range("D20").Value = "'5.256"
$ws.Range("E20").Value = "'4.44%"

# Row 21
$ws.Range("D21").Value = "'0.1269"
$ws.Range("E21").Value = "'-1.92%"

# Row 22
$ws.Range("D22").Value = "'0.2592"
$ws.Range("E22").Value = "'2.57%"

# Row 23
$ws.Range("E23").Value = "'5,592.38%"

# Row 24
$ws.Range("D24").Value = "'0.04358"
$ws.Range("E24").Value = "'-0.19%"

# Row 25
$ws.Range("D25").Value = "'0.001199"
$ws.Range("E25").Value = "'-3.07%"

# Row 26
$ws.Range("D26").Value = "'0.004251"
$ws.Range("E26").Value = "'-8.24%"

# Row 27
$ws.Range("E27").Value = "'-63.78%"

# Row 39
$ws.Range("D39").Value = "'0.02093"
$ws.Range("E39").Value = "'-5.28%"

# Row 40
$ws.Range("D40").Value = "'0.05076"
$ws.Range("E40").Value = "'0.66%"

# Row 41
$ws.Range("D41").Value = "'0.007675"
$ws.Range("E41").Value = "'1.90%"

# Row 42
$ws.Range("D42").Value = "'0.009803"
$ws.Range("E42").Value = "'-2.29%"

# Row 43
$ws.Range("D43").Value = "'0.1346"
$ws.Range("E43").Value = "'-1.41%"

# Row 44
$ws.Range("D44").Value = "'0.002082"
$ws.Range("E44").Value = "'1.98%"

# Row 45
$ws.Range("D45").Value = "'0.009589"
$ws.Range("E45").Value = "'8.46%"

# Row 46
$ws.Range("D46").Value = "'0.00006689"
$ws.Range("E46").Value = "'0.29%"

# Row 47
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.57%"

# Row 48
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001201"
$ws.Range("E48").Value = "'-0.54%"

# Row 49
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002937"
$ws.Range("E49").Value = "'-2.71%"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.57%"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.57%"

Write-Output "Applied 31-Jan-2023 09:59 UTC symbol list refresh"
